$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 259, shifting existing rows 259..343 down to 260..344
$ws.Rows(259).Insert()

# Populate the newly inserted row 259 with the new data record
$ws.Range("A259").Value = 10
$ws.Range("B259").Value = "Vega Modelo de Temuco"
$ws.Range("C259").Value = "La Araucanía"
$ws.Range("D259").Value = 44524
$ws.Range("E259").Value = 9
$ws.Range("F259").Value = 100112028
$ws.Range("G259").Value = "Sandia"
$ws.Range("H259").Value = "Sin especificar"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 250
$ws.Range("K259").Value = 1000
$ws.Range("L259").Value = 1000
$ws.Range("M259").Value = 1000
$ws.Range("N259").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O259").Value = "Región Metropolitana"
$ws.Range("P259").Value = 1000
$ws.Range("Q259").Value = 1
$ws.Range("R259").Value = "Hortaliza"

# Ensure column D keeps its existing date-number style (s="2") for the new row
$ws.Range("D259").NumberFormat = $ws.Range("D260").NumberFormat
